# datos_clientes.xlsx - re-acomodo de archivos en componentes
# Applies the row-level corrections captured in the commit diff:
#   - D2: telefono corregido
#   - J2: vencimiento corregido
#   - H3: haber pagado (0 -> 100), pasa de "deuda" (rojo) a "pago" (verde)
#   - J5: vencimiento corregido
#   - H9: haber pagado (0 -> 6500), pasa de "deuda" (rojo) a "pago" (verde)
#   - I9 / J9: fechas recalculadas
#   - K9: estado recalculado de "Vencido" a "Regular" (verde -> rojo, estilo K2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value as TEXT (matches the source file's inline-string
# cells) without leaving the cell's number format / fill behind. Excel will
# happily reinterpret digit-only or dd/mm/yyyy-looking strings as numbers or
# dates unless the cell is pre-formatted as Text ("@"); we then restore the
# original "Normal" look so no stray style sticks around.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$value)
    $range.Style = "Normal"
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# 1) D2 - Telefono
Set-TextValue $ws.Range("D2") "3537651910"

# 2) J2 - Vencimiento_Formateada
Set-TextValue $ws.Range("J2") "15/02/2024"

# 3) H3 - Haber: pasa de 0 (impago) a 100 (pagado) -> estilo verde (igual que H2)
Set-TextValue $ws.Range("H3") "100"
$ws.Range("H2").Copy()
$ws.Range("H3").PasteSpecial(-4122)  # xlPasteFormats

# 4) J5 - Vencimiento_Formateada
Set-TextValue $ws.Range("J5") "13/02/2024"

# 5) H9 - Haber: pasa de 0 (impago) a 6500 (pagado) -> estilo verde (igual que H2)
Set-TextValue $ws.Range("H9") "6500"
$ws.Range("H2").Copy()
$ws.Range("H9").PasteSpecial(-4122)  # xlPasteFormats

# 6) I9 - Fecha
Set-TextValue $ws.Range("I9") "11/02/2024"

# 7) J9 - Vencimiento_Formateada
Set-TextValue $ws.Range("J9") "11/03/2024"

# 8) K9 - Estado: pasa de "Vencido" a "Regular" -> estilo rojo (igual que K2)
$ws.Range("K9").Value = "Regular"
$ws.Range("K2").Copy()
$ws.Range("K9").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
